# Now Health Worldcare ROW benefits.xlsx -- "worldcare row rates update"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Geographical Coverage column (C2:C5): drop the combined
#     "Worldwide excluding USA/Worldwide including USA" wording in favour of
#     the single "Worldwide excluding USA" value used across all plans now.
$ws.Range("C2:C5").Value = "Worldwide excluding USA"

# --- Plan start date moved forward ---
$ws.Range("BR2").Value = "2025-04-01"

# --- Deductible (BU) / co-pay (BV) wording & rates refreshed for every plan ---
$ws.Range("BU2").Value = "NIL per year deductible (default)"

$ws.Range("BU3").Value = "USD 150 per year deductible"
$ws.Range("BV3").Value = "20% co-pay"
$ws.Range("BW3").ClearContents()

$ws.Range("BU4").Value = "USD 250 per year deductible"
$ws.Range("BV4").Value = "USD 25 excess"

$ws.Range("BU5").Value = "USD 500 per year deductible"
$ws.Range("BV5").Value = "USD 15 excess"

# --- Deductible drop-down list (column BU, Sheet1) rebuilt with the new
#     "<amount> per year deductible" wording and extra RMB/SGD tiers ---
$ws.Range("BU6").Value  = "USD 1,000 per year deductible"
$ws.Range("BU7").Value  = "USD 2,500 per year deductible"
$ws.Range("BU8").Value  = "USD 5,000 per year deductible"
$ws.Range("BU9").Value  = "USD 10,000 per year deductible"
$ws.Range("BU10").Value = "USD 15,000 per year deductible"
$ws.Range("BU11").Value = "RMB 950 per year deductible"
$ws.Range("BU12").Value = "RMB 1,570 per year deductible"
$ws.Range("BU13").Value = "RMB 3,150 per year deductible"
$ws.Range("BU14").Value = "RMB 6,300 per year deductible"
$ws.Range("BU15").Value = "RMB 15,700 per year deductible"
$ws.Range("BU16").Value = "RMB 31,500 per year deductible"
$ws.Range("BU17").Value = "RMB 63,000 per year deductible"
$ws.Range("BU18").Value = "RMB 94,500 per year deductible"
$ws.Range("BU19").Value = "RMB 195 per year deductible"
$ws.Range("BU20").Value = "RMB 325 per year deductible"
$ws.Range("BU21").Value = "RMB 650 per year deductible"
$ws.Range("BU22").Value = "SGD 1,300 per year deductible"
$ws.Range("BU23").Value = "SGD 3,250 per year deductible"
$ws.Range("BU24").Value = "SGD 6,500 per year deductible"
$ws.Range("BU25").Value = "SGD 13,000 per year deductible"
$ws.Range("BU26").Value = "SGD 19,500 per year deductible"

# --- Formatting touch-ups to match the refreshed rows ---
$ws.Range("BU2").Font.Name = "Arial"
$ws.Range("BU2").WrapText = $true

$ws.Range("BU4").Font.Name = "Arial"
$ws.Range("BU4").WrapText = $false

$ws.Range("BU5").Font.Name = "Arial"
$ws.Range("BU5").WrapText = $false

$ws.Range("BU6:BU7").Font.Name = "Arial"
$ws.Range("BU6:BU7").WrapText = $false

$ws.Range("BU1:BX1").ColumnWidth = 51

# --- View refresh: re-point the visible/active window like the author left it ---
$ws2.Activate()
$ws2.Application.ActiveWindow.Zoom = 75
$ws2.Range("A1").Select()

$ws.Activate()
$ws.Application.ActiveWindow.Zoom = 75
$ws.Range("C5").Select()
